$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing E column (n_fujian), which shifts it to H
$ws.Range("E1:G1").EntireColumn.Insert()

# Step 1: E1 header "group"
$ws.Range("E1").Value = "group"

# Step 2: column E data rows 2-29, "age<i>_vac<j>" labels, filled top-to-bottom
$colE = New-Object 'object[,]' 28,1
$colE[0,0] = "age1_vac1"
$colE[1,0] = "age1_vac2"
$colE[2,0] = "age1_vac3"
$colE[3,0] = "age1_vac4"
$colE[4,0] = "age2_vac1"
$colE[5,0] = "age2_vac2"
$colE[6,0] = "age2_vac3"
$colE[7,0] = "age2_vac4"
$colE[8,0] = "age3_vac1"
$colE[9,0] = "age3_vac2"
$colE[10,0] = "age3_vac3"
$colE[11,0] = "age3_vac4"
$colE[12,0] = "age4_vac1"
$colE[13,0] = "age4_vac2"
$colE[14,0] = "age4_vac3"
$colE[15,0] = "age4_vac4"
$colE[16,0] = "age5_vac1"
$colE[17,0] = "age5_vac2"
$colE[18,0] = "age5_vac3"
$colE[19,0] = "age5_vac4"
$colE[20,0] = "age6_vac1"
$colE[21,0] = "age6_vac2"
$colE[22,0] = "age6_vac3"
$colE[23,0] = "age6_vac4"
$colE[24,0] = "age7_vac1"
$colE[25,0] = "age7_vac2"
$colE[26,0] = "age7_vac3"
$colE[27,0] = "age7_vac4"
$ws.Range("E2:E29").Value = $colE

# Step 3: columns F and G data rows 2-29, filled row-wise (F then G each row)
$colFG = New-Object 'object[,]' 28,2
$colFG[0,0] = "age1"
$colFG[0,1] = "vac1"
$colFG[1,0] = "age1"
$colFG[1,1] = "vac2"
$colFG[2,0] = "age1"
$colFG[2,1] = "vac3"
$colFG[3,0] = "age1"
$colFG[3,1] = "vac4"
$colFG[4,0] = "age2"
$colFG[4,1] = "vac1"
$colFG[5,0] = "age2"
$colFG[5,1] = "vac2"
$colFG[6,0] = "age2"
$colFG[6,1] = "vac3"
$colFG[7,0] = "age2"
$colFG[7,1] = "vac4"
$colFG[8,0] = "age3"
$colFG[8,1] = "vac1"
$colFG[9,0] = "age3"
$colFG[9,1] = "vac2"
$colFG[10,0] = "age3"
$colFG[10,1] = "vac3"
$colFG[11,0] = "age3"
$colFG[11,1] = "vac4"
$colFG[12,0] = "age4"
$colFG[12,1] = "vac1"
$colFG[13,0] = "age4"
$colFG[13,1] = "vac2"
$colFG[14,0] = "age4"
$colFG[14,1] = "vac3"
$colFG[15,0] = "age4"
$colFG[15,1] = "vac4"
$colFG[16,0] = "age5"
$colFG[16,1] = "vac1"
$colFG[17,0] = "age5"
$colFG[17,1] = "vac2"
$colFG[18,0] = "age5"
$colFG[18,1] = "vac3"
$colFG[19,0] = "age5"
$colFG[19,1] = "vac4"
$colFG[20,0] = "age6"
$colFG[20,1] = "vac1"
$colFG[21,0] = "age6"
$colFG[21,1] = "vac2"
$colFG[22,0] = "age6"
$colFG[22,1] = "vac3"
$colFG[23,0] = "age6"
$colFG[23,1] = "vac4"
$colFG[24,0] = "age7"
$colFG[24,1] = "vac1"
$colFG[25,0] = "age7"
$colFG[25,1] = "vac2"
$colFG[26,0] = "age7"
$colFG[26,1] = "vac3"
$colFG[27,0] = "age7"
$colFG[27,1] = "vac4"
$ws.Range("F2:G29").Value = $colFG

# Step 4: headers F1/G1 typed last
$ws.Range("F1").Value = "age_group"
$ws.Range("G1").Value = "vac_group"

# Match column width of the new group columns to column B's width
$ws.Range("E1:G29").ColumnWidth = 14.5

# Leave the selection where the author ended up
$ws.Range("J7").Select()

"done"
